$wb = $excel.ActiveWorkbook

# --- 1. Update the "Status" text from "Ready for handoff" to "In Translation" ---
# This string is shared across the Overview sheet (columns E/F) and the
# per-language detail sheets (zh-cn, de-de) in column C ("Status").
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$ws1.Range("E2:F3").Value = "In Translation"
$ws2.Range("C2:C3").Value = "In Translation"
$ws3.Range("C2:C3").Value = "In Translation"

# --- 2. Narrow the affected "Status" columns to fit the new, shorter text ---
# Overview: columns E and F (zh-cn / de-de status)
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de detail sheets: column C ("Status")
$ws2.Columns.Item(3).ColumnWidth = 12.5
$ws3.Columns.Item(3).ColumnWidth = 12.5
